# Rerunning the appendix table: update the standard-error rows
# (theta_se = row 4, lambda_se = row 6) with the new bootstrap estimates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# theta_se row
$ws.Range("C4").Value = "(0.36)"
$ws.Range("D4").Value = "(0.35)"
$ws.Range("E4").Value = "(0.28)"
$ws.Range("F4").Value = "(0.41)"
$ws.Range("G4").Value = "(0.4)"

# lambda_se row
$ws.Range("C6").Value = "(0.33)"
$ws.Range("D6").Value = "(0.29)"
$ws.Range("E6").Value = "(0.27)"
$ws.Range("F6").Value = "(0.32)"
$ws.Range("G6").Value = "(0.29)"
